$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"="1.02"; "C"="1.044907436455103"; "D"="1.043104113287314"; "E"="1.058259258798868"; "F"="1.065190605937569"; "I"="1.036213688336276"; "J"="1.04997000918498"; "K"="1.045878920927564"; "L"="1.060991981814341"; "M"="1.067904513160727"; "N"="1.020377239862884" }
  3 = @{ "B"="1.02"; "C"="1.046390200558179"; "D"="1.044182475708159"; "E"="1.059770249435896"; "F"="1.066848038827471"; "I"="1.036546243454111"; "J"="1.051097584402256"; "K"="1.046767550902195"; "L"="1.062315203554813"; "M"="1.069375192515088"; "N"="1.020763360413833" }
  4 = @{ "B"="1.02"; "C"="1.04734791772583"; "D"="1.044878428622993"; "E"="1.060746614319666"; "F"="1.067919370122905"; "I"="1.036758891133731"; "J"="1.051825079152873"; "K"="1.047340146104563"; "L"="1.063169567162823"; "M"="1.07032520653673"; "N"="1.021012193540338" }
  5 = @{ "B"="1.02"; "C"="1.047750135821745"; "D"="1.04517057640113"; "E"="1.061156763300459"; "F"="1.068369493817113"; "I"="1.036847682692136"; "J"="1.052130415991772"; "K"="1.04758029324843"; "L"="1.063528305868115"; "M"="1.070724213832518"; "N"="1.021116562537635" }
  6 = @{ "B"="1.02"; "C"="1.047817646422947"; "D"="1.045219604173556"; "E"="1.061225610880315"; "F"="1.068445056334932"; "I"="1.036862555735072"; "J"="1.052181654129581"; "K"="1.047620581606921"; "L"="1.063588514324245"; "M"="1.070791186936565"; "N"="1.021134072509047" }
  7 = @{ "B"="1.02"; "C"="1.04735329376988"; "D"="1.044882334002862"; "E"="1.060752095977962"; "F"="1.067925385724389"; "I"="1.036760079947741"; "J"="1.051829161043603"; "K"="1.047343357203966"; "L"="1.063174362350587"; "M"="1.070330539569812"; "N"="1.021013589066336" }
  8 = @{ "B"="1.02"; "C"="1.045408906668425"; "D"="1.043468930599585"; "E"="1.058770188722515"; "F"="1.065750984095224"; "I"="1.03632660356311"; "J"="1.050351521476468"; "K"="1.046179738895705"; "L"="1.061439557930471"; "M"="1.068401874982852"; "N"="1.020507941997458" }
  9 = @{ "B"="1.02"; "C"="1.041969033489621"; "D"="1.040964179528591"; "E"="1.055267152062418"; "F"="1.061910303127156"; "I"="1.03554323579818"; "J"="1.047731210329054"; "K"="1.044110648099917"; "L"="1.058368132878577"; "M"="1.064990611932563"; "N"="1.019609081544482" }
  10 = @{ "B"="1.02"; "C"="1.03966612184061"; "D"="1.039284522351295"; "E"="1.052924139000165"; "F"="1.059343181922102"; "I"="1.035007732205117"; "J"="1.045972846633695"; "K"="1.042718426614907"; "L"="1.056310339760023"; "M"="1.062707400093712"; "N"="1.019004439688277" }
  11 = @{ "B"="1.02"; "C"="1.038666529711178"; "D"="1.038554813012468"; "E"="1.05190766374732"; "F"="1.058229886935084"; "I"="1.034772679016025"; "J"="1.045208644853335"; "K"="1.042112471218639"; "L"="1.05541677969784"; "M"="1.061716490729595"; "N"="1.018741313427918" }
  12 = @{ "B"="1.02"; "C"="1.038294863764729"; "D"="1.038283398808978"; "E"="1.051529799365576"; "F"="1.057816091435959"; "I"="1.034684889917484"; "J"="1.044924355325623"; "K"="1.04188691875562"; "L"="1.055084484149742"; "M"="1.061348073217804"; "N"="1.018643376957114" }
  13 = @{ "B"="1.02"; "C"="1.038374604407589"; "D"="1.038341634812957"; "E"="1.051610866296616"; "F"="1.057904864281314"; "I"="1.034703742735201"; "J"="1.044985356017109"; "K"="1.041935322051294"; "L"="1.055155780350113"; "M"="1.061427116001948"; "N"="1.018664393738014" }
  14 = @{ "B"="1.02"; "C"="1.038635815353295"; "D"="1.038532385365424"; "E"="1.051876435539082"; "F"="1.058195688002175"; "I"="1.034765432146457"; "J"="1.045185154227288"; "K"="1.042093836677505"; "L"="1.055389320000466"; "M"="1.061686044411679"; "N"="1.01873322205626" }
  15 = @{ "B"="1.02"; "C"="1.038796706323928"; "D"="1.038649864174053"; "E"="1.052040021447843"; "F"="1.05837483824166"; "I"="1.034803377360964"; "J"="1.04530819916726"; "K"="1.042191439833987"; "L"="1.055533159872071"; "M"="1.061845532130141"; "N"="1.018775602922509" }
  16 = @{ "B"="1.02"; "C"="1.039732409727296"; "D"="1.039332899532548"; "E"="1.052991557536515"; "F"="1.059417030619178"; "I"="1.035023264760075"; "J"="1.046023504165185"; "K"="1.042758575773806"; "L"="1.056369588525502"; "M"="1.062773114978304"; "N"="1.019021874649371" }
  17 = @{ "B"="1.02"; "C"="1.040318697894181"; "D"="1.039760700949658"; "E"="1.053587906226938"; "F"="1.060070303886472"; "I"="1.035160342003308"; "J"="1.046471436185752"; "K"="1.043113487329589"; "L"="1.056893576815667"; "M"="1.06335435043481"; "N"="1.019176001289233" }
  18 = @{ "B"="1.02"; "C"="1.040660437410436"; "D"="1.040009998482214"; "E"="1.053935560261533"; "F"="1.060451182419854"; "I"="1.035239990524099"; "J"="1.046732435901169"; "K"="1.043320201142822"; "L"="1.057198967566122"; "M"="1.063693157492927"; "N"="1.019265774298317" }
  19 = @{ "B"="1.02"; "C"="1.040776922624693"; "D"="1.040094963364566"; "E"="1.054054070028262"; "F"="1.060581024537489"; "I"="1.035267096733789"; "J"="1.04682138424111"; "K"="1.043390634483334"; "L"="1.057303057004803"; "M"="1.063808645298246"; "N"="1.019296363201045" }
  20 = @{ "B"="1.02"; "C"="1.040255818837883"; "D"="1.039714825926133"; "E"="1.053523942970699"; "F"="1.060000230960817"; "I"="1.035145666615149"; "J"="1.046423405451055"; "K"="1.043075439760202"; "L"="1.056837383002885"; "M"="1.063292011963331"; "N"="1.019159478057064" }
  21 = @{ "B"="1.02"; "C"="1.038558905614869"; "D"="1.038476224280593"; "E"="1.051798240380267"; "F"="1.058110055153221"; "I"="1.034747279420682"; "J"="1.045126330591287"; "K"="1.04204717119253"; "L"="1.055320559162487"; "M"="1.061609806178127"; "N"="1.018712959363852" }
  22 = @{ "B"="1.02"; "C"="1.037489826329846"; "D"="1.037695336010245"; "E"="1.050711480963636"; "F"="1.056920069738679"; "I"="1.034494019825314"; "J"="1.044308310985135"; "K"="1.041397914650887"; "L"="1.054364625198199"; "M"="1.06055010884031"; "N"="1.018431058902831" }
  23 = @{ "B"="1.02"; "C"="1.038056773808173"; "D"="1.038109503628765"; "E"="1.051287760972971"; "F"="1.057551054649403"; "I"="1.034628541713795"; "J"="1.044742197873626"; "K"="1.041742359797142"; "L"="1.054871599737995"; "M"="1.061112069879621"; "N"="1.018580610127112" }
  24 = @{ "B"="1.02"; "C"="1.040284231895262"; "D"="1.039735555591069"; "E"="1.053552845788069"; "F"="1.06003189440881"; "I"="1.035152298737964"; "J"="1.046445109311552"; "K"="1.043092632746515"; "L"="1.056862775318326"; "M"="1.063320180708285"; "N"="1.019166944585253" }
  25 = @{ "B"="1.02"; "C"="1.042859988724301"; "D"="1.041613426407595"; "E"="1.056174084410837"; "F"="1.062904346221509"; "I"="1.035748082426176"; "J"="1.04841062050915"; "K"="1.044647797128163"; "L"="1.059163930886444"; "M"="1.065874061674931"; "N"="1.019842401231313" }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$row").Value = [double]$rowData[$col]
  }
}
